# Update the cryptocurrency price/volume table on Sheet1 to reflect
# refreshed market data (GitHub Actions scheduled update).
#
# Price values in column D are stored as text (they use a dotted
# thousands/decimal style like "36.170.35" and must not be reinterpreted
# as numbers), so they are written with a leading apostrophe to force
# Excel to keep them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''36.170.35'
$ws.Range("E2").Value = '  -3.61%  '
$ws.Range("D3").Value = '''1.950.40'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''229.18'
$ws.Range("E5").Value = '  -9.84%  '
$ws.Range("D6").Value = '''0.590'
$ws.Range("E6").Value = '  -4.77%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = '''52.81'
$ws.Range("E8").Value = '  -6.91%  '
$ws.Range("D9").Value = '''0.363'
$ws.Range("E9").Value = '  -4.75%  '
$ws.Range("D10").Value = '''56.47'
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").Value = '''0.0729'
$ws.Range("E11").Value = '  -6.82%  '
$ws.Range("D12").Value = '''0.0968'
$ws.Range("E12").Value = '  -4.63%  '
$ws.Range("D13").Value = '''2.235.70'
$ws.Range("E13").Value = '  -3.71%  '
$ws.Range("D14").Value = '''13.69'
$ws.Range("E14").Value = '  -5.57%  '
$ws.Range("D15").Value = '''19.31'
$ws.Range("E15").Value = '  -8.12%  '
$ws.Range("D16").Value = '''0.736'
$ws.Range("E16").Value = '  -9.57%  '
$ws.Range("D17").Value = '''1.958.99'
$ws.Range("E17").Value = '  -3.58%  '
$ws.Range("D18").Value = '''4.93'
$ws.Range("E18").Value = '  -7.56%  '
$ws.Range("D19").Value = '''36.143.27'
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").Value = '''66.60'
$ws.Range("E20").Value = '  -4.24%  '
$ws.Range("D21").Value = '''0.0₃0783'
$ws.Range("E21").Value = '  -7.60%  '
$ws.Range("D22").Value = '''4.92'
$ws.Range("E22").Value = '  -5.20%  '
$ws.Range("D23").Value = '''218.90'
$ws.Range("E23").Value = '  -4.12%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").Value = '''2.30'
$ws.Range("E26").Value = '  -12.00%  '
$ws.Range("D27").Value = '''159.82'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").Value = '''8.36'
$ws.Range("E28").Value = '  -7.39%  '
$ws.Range("D29").Value = '''18.58'
$ws.Range("E29").Value = '  -6.09%  '
$ws.Range("E30").Value = '  -7.57%  '
$ws.Range("E31").Value = '  -11.34%  '
$ws.Range("E32").Value = '  -4.15%  '
$ws.Range("E33").Value = '  -9.34%  '
$ws.Range("D34").Value = '''0.0596'
$ws.Range("E34").Value = '  -10.57%  '
$ws.Range("D35").Value = '''4.15'
$ws.Range("E35").Value = '  -8.97%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '''2.23'
$ws.Range("E36").Value = '  -7.99%  '
$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.20%  '
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("D39").Value = '''3.12'
$ws.Range("E39").Value = '  -7.86%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").Value = '''4.98'
$ws.Range("E40").Value = '  -6.94%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").Value = '''2.99'
$ws.Range("E41").Value = '  -1.49%  '
$ws.Range("D42").Value = '''1.395.81'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '''0.0197'
$ws.Range("E43").Value = '  -8.55%  '
$ws.Range("E44").Value = '  -11.51%  '
$ws.Range("E45").Value = '  -11.21%  '
$ws.Range("D46").Value = '''85.55'
$ws.Range("E46").Value = '  -5.64%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '''0.966'
$ws.Range("E47").Value = '  -6.61%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '''2.85'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = '''14.50'
$ws.Range("E49").Value = '  -9.03%  '
$ws.Range("D50").Value = '''6.66'
$ws.Range("E50").Value = '  -8.82%  '
$ws.Range("D51").Value = '''2.129.78'
$ws.Range("E51").Value = '  -3.81%  '
